$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) contain text values that look numeric (e.g. "213.68" or "26.004.45").
# Excel auto-converts such strings to numbers on assignment, so force the cell to Text
# format first for every Price cell we touch, to preserve the original string data type.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.997.77"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.627.62"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "213.68"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").Value = "0.250"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("E10").Value = "  -5.24%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "1.854.26"
$ws.Range("D13").Value = "1.632.33"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "25.998.99"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "61.45"
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "192.39"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "144.23"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("D29").Value = "15.22"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  -4.79%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").Value = "1.123.65"
$ws.Range("E37").Value = "  -5.56%  "
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").Value = "98.33"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "1.764.45"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  -4.50%  "
$ws.Range("E44").Value = "  -5.83%  "
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "54.36"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.48"
$ws.Range("E51").Value = "  -3.55%  "
